$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 229.3125
$ws.Range("I33").Value = 223.6
$ws.Range("J33").Value = 238.83333
$ws.Range("K33").Value = 223.6
$ws.Range("L33").Value = 238.83333
$ws.Range("M33").Value = 5.400000000000006
$ws.Range("N33").Value = -696.8333299999999

$ws.Range("H70").Value = 92018770
$ws.Range("I70").Value = 337398600
$ws.Range("J70").Value = 1325
$ws.Range("K70").Value = 1012195800
$ws.Range("L70").Value = 3975
$ws.Range("M70").Value = -1012195530

$ws.Range("H73").Value = 92018770
$ws.Range("I73").Value = 337398600
$ws.Range("J73").Value = 1325
$ws.Range("K73").Value = 1012195800
$ws.Range("L73").Value = 3975
$ws.Range("M73").Value = -1012194864

$ws.Range("H112").Value = 5555.143
$ws.Range("I112").Value = 4929.3335
$ws.Range("J112").Value = 6024.5
$ws.Range("K112").Value = 14788.0005
$ws.Range("L112").Value = 18073.5
$ws.Range("M112").Value = -13680.0005
$ws.Range("N112").Value = -20289.5

$ws.Range("H121").Value = 490
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 490
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 1470
$ws.Range("N121").Value = -4964

$ws.Range("H137").Value = 3988.3125
$ws.Range("I137").Value = 2897.348
$ws.Range("J137").Value = 6776.3335
$ws.Range("K137").Value = 8692.044
$ws.Range("L137").Value = 20329.0005
$ws.Range("M137").Value = -6142.044
$ws.Range("N137").Value = -25429.0005

$ws.Range("H138").Value = 7293.0513
$ws.Range("I138").Value = 5412.5884
$ws.Range("J138").Value = 8746.137000000001
$ws.Range("K138").Value = 16237.7652
$ws.Range("L138").Value = 26238.411
$ws.Range("M138").Value = -11097.7652
$ws.Range("N138").Value = -36518.411

$ws.Range("H141").Value = 13528186
$ws.Range("I141").Value = 14709545
$ws.Range("J141").Value = 139444
$ws.Range("K141").Value = 44128635
$ws.Range("L141").Value = 418332
$ws.Range("M141").Value = -44123455
$ws.Range("N141").Value = -428692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1385.25
$ws.Range("I5").Value = 1680.3334
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 1680.3334
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -1568.3334
$ws.Range("N5").Value = -724

$ws.Range("H32").Value = 2989.3623
$ws.Range("I32").Value = 2913.1365
$ws.Range("J32").Value = 4666.3335
$ws.Range("K32").Value = 2913.1365
$ws.Range("L32").Value = 4666.3335
$ws.Range("M32").Value = -2626.1365
$ws.Range("N32").Value = -5240.3335

$ws.Range("H45").Value = 1834983.9
$ws.Range("I45").Value = 3403201.2
$ws.Range("J45").Value = 5396.8335
$ws.Range("K45").Value = 3403201.2
$ws.Range("L45").Value = 5396.8335
$ws.Range("M45").Value = -3402824.2
$ws.Range("N45").Value = -6150.8335

$ws.Range("H74").Value = 1887.561
$ws.Range("I74").Value = 1718.0312
$ws.Range("J74").Value = 2490.3333
$ws.Range("K74").Value = 1718.0312
$ws.Range("L74").Value = 2490.3333
$ws.Range("M74").Value = -844.0311999999999
$ws.Range("N74").Value = -4238.3333

$ws.Range("H77").Value = 1887.561
$ws.Range("I77").Value = 1718.0312
$ws.Range("J77").Value = 2490.3333
$ws.Range("K77").Value = 8590.155999999999
$ws.Range("L77").Value = 12451.6665
$ws.Range("M77").Value = -4222.155999999999
$ws.Range("N77").Value = -21187.6665

$ws.Range("H122").Value = 3878.675
$ws.Range("I122").Value = 3882.8438
$ws.Range("J122").Value = 3862
$ws.Range("K122").Value = 11648.5314
$ws.Range("L122").Value = 11586
$ws.Range("M122").Value = -9198.5314

$ws.Range("H132").Value = 1964525.2
$ws.Range("I132").Value = 3580.617
$ws.Range("J132").Value = 25005626
$ws.Range("K132").Value = 10741.851
$ws.Range("L132").Value = 75016878
$ws.Range("M132").Value = -8211.851000000001
$ws.Range("N132").Value = -75021938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1385.25
$ws.Range("I4").Value = 1680.3334
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1680.3334
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -1565.3334
$ws.Range("N4").Value = -730

$ws.Range("H22").Value = 116
$ws.Range("I22").Value = 82.5
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 82.5
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 90.5

$ws.Range("H26").Value = 64443.5
$ws.Range("I26").Value = 28888
$ws.Range("J26").Value = 99999
$ws.Range("K26").Value = 28888
$ws.Range("L26").Value = 99999
$ws.Range("M26").Value = -28596

$ws.Range("H132").Value = 195000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 195000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 195000
$ws.Range("N132").Value = -205120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21758250
$ws.Range("I31").Value = 41700300
$ws.Range("J31").Value = 3283.1365
$ws.Range("K31").Value = 41700300
$ws.Range("L31").Value = 3283.1365
$ws.Range("M31").Value = -41700005
$ws.Range("N31").Value = -3873.1365

$ws.Range("H34").Value = 21758250
$ws.Range("I34").Value = 41700300
$ws.Range("J34").Value = 3283.1365
$ws.Range("K34").Value = 41700300
$ws.Range("L34").Value = 3283.1365
$ws.Range("M34").Value = -41700098
$ws.Range("N34").Value = -3687.1365

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1277.6
$ws.Range("I50").Value = 847.5
$ws.Range("J50").Value = 2998
$ws.Range("K50").Value = 2542.5
$ws.Range("L50").Value = 8994
$ws.Range("M50").Value = -2061.5
$ws.Range("N50").Value = -9956

$ws.Range("H53").Value = 1277.6
$ws.Range("I53").Value = 847.5
$ws.Range("J53").Value = 2998
$ws.Range("K53").Value = 2542.5
$ws.Range("L53").Value = 8994
$ws.Range("M53").Value = -2061.5
$ws.Range("N53").Value = -9956

$ws.Range("H56").Value = 15097.3125
$ws.Range("I56").Value = 15097.3125
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 15097.3125
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -14567.3125

$ws.Range("H68").Value = 3006.6924
$ws.Range("I68").Value = 1700
$ws.Range("J68").Value = 3115.5833
$ws.Range("K68").Value = 5100
$ws.Range("L68").Value = 9346.749899999999
$ws.Range("M68").Value = -4289
$ws.Range("N68").Value = -10968.7499

$ws.Range("H71").Value = 3006.6924
$ws.Range("I71").Value = 1700
$ws.Range("J71").Value = 3115.5833
$ws.Range("K71").Value = 15300
$ws.Range("L71").Value = 28040.2497
$ws.Range("M71").Value = -11244
$ws.Range("N71").Value = -36152.2497

$ws.Range("H112").Value = 14731
$ws.Range("I112").Value = 8351
$ws.Range("J112").Value = 21111
$ws.Range("K112").Value = 25053
$ws.Range("L112").Value = 63333
$ws.Range("M112").Value = -23945

$ws.Range("H137").Value = 20261.715
$ws.Range("I137").Value = 7499.5
$ws.Range("J137").Value = 25366.6
$ws.Range("K137").Value = 22498.5
$ws.Range("L137").Value = 76099.79999999999
$ws.Range("M137").Value = -17398.5
$ws.Range("N137").Value = -86299.79999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1142.5652
$ws.Range("I97").Value = 950.2
$ws.Range("J97").Value = 1503.25
$ws.Range("K97").Value = 950.2
$ws.Range("L97").Value = 1503.25
$ws.Range("M97").Value = -454.2
$ws.Range("N97").Value = -2495.25

$ws.Range("H132").Value = 3410768.2
$ws.Range("I132").Value = 1680.2444
$ws.Range("J132").Value = 17357038
$ws.Range("K132").Value = 5040.733200000001
$ws.Range("L132").Value = 52071114
$ws.Range("M132").Value = -2510.733200000001
$ws.Range("N132").Value = -52076174

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1465.7142
$ws.Range("I46").Value = 1465.7142
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1465.7142
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1277.7142
$ws.Range("N46").ClearContents()

$ws.Range("H55").Value = 601.0833
$ws.Range("I55").Value = 372.04166
$ws.Range("J55").Value = 830.125
$ws.Range("K55").Value = 372.04166
$ws.Range("L55").Value = 830.125
$ws.Range("M55").Value = -199.04166
$ws.Range("N55").Value = -1176.125

$ws.Range("H68").Value = 2195111.5
$ws.Range("I68").Value = 4168258.2
$ws.Range("J68").Value = 2726.2222
$ws.Range("K68").Value = 4168258.2
$ws.Range("L68").Value = 2726.2222
$ws.Range("M68").Value = -4167509.2
$ws.Range("N68").Value = -4224.2222

$ws.Range("H71").Value = 2195111.5
$ws.Range("I71").Value = 4168258.2
$ws.Range("J71").Value = 2726.2222
$ws.Range("K71").Value = 20841291
$ws.Range("L71").Value = 13631.111
$ws.Range("M71").Value = -20837547
$ws.Range("N71").Value = -21119.111

$ws.Range("H136").Value = 7489.3125
$ws.Range("I136").Value = 8035.0713
$ws.Range("J136").Value = 3669
$ws.Range("K136").Value = 24105.2139
$ws.Range("L136").Value = 11007
$ws.Range("M136").Value = -21555.2139
